$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.965.48"
$ws.Range("E2").Value = "  -4.76%  "
$ws.Range("D3").Value = "2.226.92"
$ws.Range("E3").Value = "  -5.71%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.593"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.90%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0831"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.29%  "
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.864"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -11.39%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.564.83"
$ws.Range("E16").Value = "  -5.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.75%  "
$ws.Range("D18").Value = "2.223.52"
$ws.Range("E18").Value = "  -5.73%  "
$ws.Range("D19").Value = "42.852.44"
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.17%  "
$ws.Range("D21").Value = "0.0₃0965"
$ws.Range("E21").Value = "  -8.67%  "
$ws.Range("E22").Value = "  -10.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -9.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "236.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.54%  "
$ws.Range("E26").Value = "  -6.77%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.02%  "
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0906"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.91%  "
$ws.Range("E36").Value = "  +10.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.61%  "
$ws.Range("E38").Value = "  -5.74%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.105"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.21%  "
$ws.Range("E42").Value = "  -6.92%  "
$ws.Range("D43").Value = "1.936.57"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.209"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.875"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +19.79%  "
